$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2,1).Value = "Última actualización: 08:55:25"
$ws1.Cells.Item(3,1).Value = "Total filas: 111"

$ws1.Cells.Item(79,1).Value = "08:55:25"
$ws1.Cells.Item(79,2).Value = "08:55"
$ws1.Cells.Item(79,3).Value = "10_OLMOS"
$ws1.Cells.Item(79,4).Value = 0
$ws1.Cells.Item(79,5).Value = "LP1912"
$ws1.Cells.Item(80,1).Value = "08:55:25"
$ws1.Cells.Item(80,2).Value = "08:56"
$ws1.Cells.Item(80,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(80,4).Value = 1
$ws1.Cells.Item(80,5).Value = "LP1912"
$ws1.Cells.Item(81,1).Value = "07:24:45"
$ws1.Cells.Item(81,2).Value = "09:01"
$ws1.Cells.Item(81,3).Value = "215A_EL PATO"
$ws1.Cells.Item(81,4).Value = 97
$ws1.Cells.Item(81,5).Value = "LP1912"
$ws1.Cells.Item(82,1).Value = "08:00:32"
$ws1.Cells.Item(82,2).Value = "09:03"
$ws1.Cells.Item(82,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(82,4).Value = 63
$ws1.Cells.Item(82,5).Value = "LP1912"
$ws1.Cells.Item(83,1).Value = "08:31:53"
$ws1.Cells.Item(83,2).Value = "09:04"
$ws1.Cells.Item(83,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(83,4).Value = 33
$ws1.Cells.Item(83,5).Value = "LP1912"
$ws1.Cells.Item(84,1).Value = "08:47:51"
$ws1.Cells.Item(84,2).Value = "09:05"
$ws1.Cells.Item(84,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(84,4).Value = 18
$ws1.Cells.Item(84,5).Value = "LP1912"
$ws1.Cells.Item(85,1).Value = "08:55:25"
$ws1.Cells.Item(85,2).Value = "09:06"
$ws1.Cells.Item(85,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(85,4).Value = 11
$ws1.Cells.Item(85,5).Value = "LP1912"
$ws1.Cells.Item(86,1).Value = "07:48:31"
$ws1.Cells.Item(86,2).Value = "09:07"
$ws1.Cells.Item(86,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(86,4).Value = 79
$ws1.Cells.Item(86,5).Value = "LP1912"
$ws1.Cells.Item(87,1).Value = "08:00:32"
$ws1.Cells.Item(87,2).Value = "09:08"
$ws1.Cells.Item(87,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(87,4).Value = 68
$ws1.Cells.Item(87,5).Value = "LP1912"
$ws1.Cells.Item(88,1).Value = "07:24:45"
$ws1.Cells.Item(88,2).Value = "09:10"
$ws1.Cells.Item(88,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(88,4).Value = 106
$ws1.Cells.Item(88,5).Value = "LP1912"
$ws1.Cells.Item(89,1).Value = "08:47:51"
$ws1.Cells.Item(89,2).Value = "09:13"
$ws1.Cells.Item(89,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(89,4).Value = 26
$ws1.Cells.Item(89,5).Value = "LP1912"
$ws1.Cells.Item(90,1).Value = "07:24:45"
$ws1.Cells.Item(90,2).Value = "09:16"
$ws1.Cells.Item(90,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(90,4).Value = 112
$ws1.Cells.Item(90,5).Value = "LP1912"
$ws1.Cells.Item(91,1).Value = "08:55:25"
$ws1.Cells.Item(91,2).Value = "09:16"
$ws1.Cells.Item(91,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(91,4).Value = 21
$ws1.Cells.Item(91,5).Value = "LP1912"
$ws1.Cells.Item(92,1).Value = "08:55:25"
$ws1.Cells.Item(92,2).Value = "09:17"
$ws1.Cells.Item(92,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(92,4).Value = 22
$ws1.Cells.Item(92,5).Value = "LP1912"
$ws1.Cells.Item(93,1).Value = "08:31:53"
$ws1.Cells.Item(93,2).Value = "09:20"
$ws1.Cells.Item(93,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(93,4).Value = 49
$ws1.Cells.Item(93,5).Value = "LP1912"
$ws1.Cells.Item(94,1).Value = "07:24:45"
$ws1.Cells.Item(94,2).Value = "09:21"
$ws1.Cells.Item(94,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(94,4).Value = 117
$ws1.Cells.Item(94,5).Value = "LP1912"
$ws1.Cells.Item(95,1).Value = "08:00:32"
$ws1.Cells.Item(95,2).Value = "09:22"
$ws1.Cells.Item(95,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(95,4).Value = 82
$ws1.Cells.Item(95,5).Value = "LP1912"
$ws1.Cells.Item(96,1).Value = "07:24:45"
$ws1.Cells.Item(96,2).Value = "09:22"
$ws1.Cells.Item(96,3).Value = "17_ROMERO"
$ws1.Cells.Item(96,4).Value = 118
$ws1.Cells.Item(96,5).Value = "LP1912"
$ws1.Cells.Item(97,1).Value = "07:48:31"
$ws1.Cells.Item(97,2).Value = "09:23"
$ws1.Cells.Item(97,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(97,4).Value = 95
$ws1.Cells.Item(97,5).Value = "LP1912"
$ws1.Cells.Item(98,1).Value = "07:48:31"
$ws1.Cells.Item(98,2).Value = "09:32"
$ws1.Cells.Item(98,3).Value = "15_ABASTO"
$ws1.Cells.Item(98,4).Value = 104
$ws1.Cells.Item(98,5).Value = "LP1912"
$ws1.Cells.Item(99,1).Value = "07:48:31"
$ws1.Cells.Item(99,2).Value = "09:33"
$ws1.Cells.Item(99,3).Value = "10_OLMOS"
$ws1.Cells.Item(99,4).Value = 105
$ws1.Cells.Item(99,5).Value = "LP1912"
$ws1.Cells.Item(100,1).Value = "08:47:51"
$ws1.Cells.Item(100,2).Value = "09:34"
$ws1.Cells.Item(100,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(100,4).Value = 47
$ws1.Cells.Item(100,5).Value = "LP1912"
$ws1.Cells.Item(101,1).Value = "08:31:53"
$ws1.Cells.Item(101,2).Value = "09:41"
$ws1.Cells.Item(101,3).Value = "215C_EL PATO"
$ws1.Cells.Item(101,4).Value = 70
$ws1.Cells.Item(101,5).Value = "LP1912"
$ws1.Cells.Item(102,1).Value = "07:48:31"
$ws1.Cells.Item(102,2).Value = "09:42"
$ws1.Cells.Item(102,3).Value = "215C_EL PATO"
$ws1.Cells.Item(102,4).Value = 114
$ws1.Cells.Item(102,5).Value = "LP1912"
$ws1.Cells.Item(103,1).Value = "08:00:32"
$ws1.Cells.Item(103,2).Value = "09:43"
$ws1.Cells.Item(103,3).Value = "14_ABASTO"
$ws1.Cells.Item(103,4).Value = 103
$ws1.Cells.Item(103,5).Value = "LP1912"
$ws1.Cells.Item(104,1).Value = "08:31:53"
$ws1.Cells.Item(104,2).Value = "09:46"
$ws1.Cells.Item(104,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(104,4).Value = 75
$ws1.Cells.Item(104,5).Value = "LP1912"
$ws1.Cells.Item(105,1).Value = "08:47:51"
$ws1.Cells.Item(105,2).Value = "09:52"
$ws1.Cells.Item(105,3).Value = "15_ABASTO"
$ws1.Cells.Item(105,4).Value = 65
$ws1.Cells.Item(105,5).Value = "LP1912"
$ws1.Cells.Item(106,1).Value = "08:47:51"
$ws1.Cells.Item(106,2).Value = "09:53"
$ws1.Cells.Item(106,3).Value = "10_OLMOS"
$ws1.Cells.Item(106,4).Value = 66
$ws1.Cells.Item(106,5).Value = "LP1912"
$ws1.Cells.Item(107,1).Value = "08:31:53"
$ws1.Cells.Item(107,2).Value = "10:03"
$ws1.Cells.Item(107,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(107,4).Value = 92
$ws1.Cells.Item(107,5).Value = "LP1912"
$ws1.Cells.Item(108,1).Value = "08:31:53"
$ws1.Cells.Item(108,2).Value = "10:10"
$ws1.Cells.Item(108,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(108,4).Value = 99
$ws1.Cells.Item(108,5).Value = "LP1912"
$ws1.Cells.Item(109,1).Value = "08:31:53"
$ws1.Cells.Item(109,2).Value = "10:12"
$ws1.Cells.Item(109,3).Value = "15_ABASTO"
$ws1.Cells.Item(109,4).Value = 101
$ws1.Cells.Item(109,5).Value = "LP1912"
$ws1.Cells.Item(110,1).Value = "08:31:53"
$ws1.Cells.Item(110,2).Value = "10:20"
$ws1.Cells.Item(110,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(110,4).Value = 109
$ws1.Cells.Item(110,5).Value = "LP1912"
$ws1.Cells.Item(111,1).Value = "08:55:25"
$ws1.Cells.Item(111,2).Value = "10:21"
$ws1.Cells.Item(111,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(111,4).Value = 86
$ws1.Cells.Item(111,5).Value = "LP1912"
$ws1.Cells.Item(112,1).Value = "08:31:53"
$ws1.Cells.Item(112,2).Value = "10:22"
$ws1.Cells.Item(112,3).Value = "17_ROMERO"
$ws1.Cells.Item(112,4).Value = 111
$ws1.Cells.Item(112,5).Value = "LP1912"
$ws1.Cells.Item(113,1).Value = "08:31:53"
$ws1.Cells.Item(113,2).Value = "10:26"
$ws1.Cells.Item(113,3).Value = "215A_EL PATO"
$ws1.Cells.Item(113,4).Value = 115
$ws1.Cells.Item(113,5).Value = "LP1912"
$ws1.Cells.Item(114,1).Value = "08:47:51"
$ws1.Cells.Item(114,2).Value = "10:41"
$ws1.Cells.Item(114,3).Value = "17_ROMERO"
$ws1.Cells.Item(114,4).Value = 114
$ws1.Cells.Item(114,5).Value = "LP1912"
$ws1.Cells.Item(115,1).Value = "08:55:25"
$ws1.Cells.Item(115,2).Value = "10:42"
$ws1.Cells.Item(115,3).Value = "17_ROMERO"
$ws1.Cells.Item(115,4).Value = 107
$ws1.Cells.Item(115,5).Value = "LP1912"
$ws1.Cells.Item(116,1).Value = "08:47:51"
$ws1.Cells.Item(116,2).Value = "10:43"
$ws1.Cells.Item(116,3).Value = "14_ABASTO"
$ws1.Cells.Item(116,4).Value = 116
$ws1.Cells.Item(116,5).Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2,1).Value = "Última actualización: 08:55:25"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2,1).Value = "Última actualización: 08:55:25"
$ws3.Cells.Item(3,1).Value = "Total filas: 26"

$ws3.Cells.Item(28,1).Value = "08:55:25"
$ws3.Cells.Item(28,2).Value = "09:10"
$ws3.Cells.Item(28,3).Value = "215D_LA PLATA"
$ws3.Cells.Item(28,4).Value = 15
$ws3.Cells.Item(28,5).Value = "L6203"
$ws3.Cells.Item(29,1).Value = "08:31:53"
$ws3.Cells.Item(29,2).Value = "10:02"
$ws3.Cells.Item(29,3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(29,4).Value = 91
$ws3.Cells.Item(29,5).Value = "L6173"
$ws3.Cells.Item(30,1).Value = "08:55:25"
$ws3.Cells.Item(30,2).Value = "10:03"
$ws3.Cells.Item(30,3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(30,4).Value = 68
$ws3.Cells.Item(30,5).Value = "L6173"
$ws3.Cells.Item(31,1).Value = "08:55:25"
$ws3.Cells.Item(31,2).Value = "10:54"
$ws3.Cells.Item(31,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(31,4).Value = 119
$ws3.Cells.Item(31,5).Value = "L6173"
